{"js": "// Fix the mis-decoded \"chi\" (\u03c7) character in the \u03c7\u00b2 column headers and\n// shrink the matching header row height from 637 (31.85pt) to 571 (28.55pt)\n// twips, for both ranova tables that contain a \"\u03c72\" header cell.\n\n// 1) Replace the mojibake \"\u00cf\u2021\" with the correct Greek small letter chi \"\u03c7\".\nconst mojibake = \"\\u00CF\\u2021\"; // \"\u00cf\u2021\"\nconst chi = \"\\u03C7\"; // \"\u03c7\"\nconst results = context.document.body.search(mojibake, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(chi, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Shrink the header row of every table whose first row's preferred\n// height is currently 31.85pt (637 twips) down to 28.55pt (571 twips).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst targetOldHeight = 31.85; // 637 twips\nconst targetNewHeight = 28.55; // 571 twips\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const headerRow = tables.items[i].rows.getFirst();\n  headerRow.load(\"preferredHeight\");\n  await context.sync();\n\n  if (Math.abs(headerRow.preferredHeight - targetOldHeight) < 0.01) {\n    headerRow.preferredHeight = targetNewHeight;\n  }\n}\nawait context.sync();\n", "ps1": "# Fix the mis-decoded \"chi\" (\u03c7) character in the \u03c7\u00b2 column headers and\n# shrink the matching header row height from 637 (31.85pt) to 571 (28.55pt)\n# twips, for both ranova tables that contain a \"\u03c72\" header cell.\n\n$d = $word.ActiveDocument\n\n# 1) Replace the mojibake \"\u00cf\u2021\" with the correct Greek small letter chi \"\u03c7\"\n#    everywhere it occurs in the document body.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u00cf\u2021\"\n$find.Replacement.Text = \"\u03c7\"\n$find.Forward = $true\n$find.Wrap = 1\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n# 2) Shrink the header row of every table whose first row's preferred\n#    height is currently 31.85pt (637 twips) down to 28.55pt (571 twips).\n$targetOldHeight = 31.85\n$targetNewHeight = 28.55\n\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $table = $d.Tables.Item($i)\n    $headerRow = $table.Rows.Item(1)\n    if ([Math]::Abs($headerRow.Height - $targetOldHeight) -lt 0.01) {\n        $headerRow.Height = $targetNewHeight\n    }\n}\n"}
